$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.87"
$ws.Range("E2").Value = "  +3.59%  "
$ws.Range("D3").Value = "1.605.49"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").Value = "1.829.41"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "1.602.00"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "26.196.07"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "197.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.33%  "
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.129"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").Value = "1.106.91"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.498"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.774"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "1.741.25"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  +10.52%  "
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -0.35%  "
